$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.625.46"
$ws.Range("E2").Value = "  +4.65%  "
$ws.Range("D3").Value = "3.247.02"
$ws.Range("E3").Value = "  +0.25%  "
$ws.Range("D4").Value = "'1.01"
$ws.Range("E4").Value = "  +0.68%  "
$ws.Range("D5").Value = "'402.67"
$ws.Range("E5").Value = "  +2.08%  "
$ws.Range("D6").Value = "'109.24"
$ws.Range("E6").Value = "  +1.67%  "
$ws.Range("D7").Value = "3.245.90"
$ws.Range("E7").Value = "  +0.38%  "
$ws.Range("D8").Value = "'0.554"
$ws.Range("E8").Value = "  -3.31%  "
$ws.Range("E9").Value = "  +0.21%  "
$ws.Range("D10").Value = "'0.602"
$ws.Range("E10").Value = "  -2.26%  "
$ws.Range("D11").Value = "'0.104"
$ws.Range("E11").Value = "  +8.41%  "
$ws.Range("D12").Value = "'37.67"
$ws.Range("E12").Value = "  -3.58%  "
$ws.Range("E13").Value = "  +0.11%  "
$ws.Range("D14").Value = "3.840.18"
$ws.Range("E14").Value = "  +2.57%  "
$ws.Range("D15").Value = "'7.96"
$ws.Range("E15").Value = "  -2.81%  "
$ws.Range("D16").Value = "'18.53"
$ws.Range("E16").Value = "  -2.92%  "
$ws.Range("D17").Value = "3.290.21"
$ws.Range("E17").Value = "  +2.08%  "
$ws.Range("D18").Value = "59.931.75"
$ws.Range("E18").Value = "  +5.46%  "
$ws.Range("D19").Value = "'0.970"
$ws.Range("E19").Value = "  -6.19%  "
$ws.Range("D20").Value = "'10.27"
$ws.Range("E20").Value = "  -6.15%  "
$ws.Range("D21").Value = "'0.0000107"
$ws.Range("E21").Value = "  +2.33%  "
$ws.Range("D22").Value = "'3.15"
$ws.Range("E22").Value = "  -5.60%  "
$ws.Range("D23").Value = "'289.70"
$ws.Range("E23").Value = "  -2.44%  "
$ws.Range("D24").Value = "'11.92"
$ws.Range("E24").Value = "  -8.10%  "
$ws.Range("D25").Value = "'72.18"
$ws.Range("E25").Value = "  -2.44%  "
$ws.Range("D26").Value = "'3.02"
$ws.Range("E26").Value = "  -4.38%  "
$ws.Range("D27").Value = "'4.48"
$ws.Range("E27").Value = "  +2.79%  "
$ws.Range("D28").Value = "'27.83"
$ws.Range("E28").Value = "  -0.02%  "
$ws.Range("D29").Value = "'7.18"
$ws.Range("E29").Value = "  -1.04%  "
$ws.Range("B30").Value = "Filecoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D30").Value = "'7.32"
$ws.Range("E30").Value = "  -4.79%  "
$ws.Range("B31").Value = "Kaspa"
$ws.Range("C31").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D31").Value = "'0.165"
$ws.Range("E31").Value = "  -1.96%  "
$ws.Range("D32").Value = "'0.999"
$ws.Range("E32").Value = "  -0.04%  "
$ws.Range("B33").Value = "Cosmos"
$ws.Range("C33").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D33").Value = "'10.91"
$ws.Range("E33").Value = "  -4.40%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").Value = "'0.106"
$ws.Range("E34").Value = "  -2.17%  "
$ws.Range("B35").Value = "Toncoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D35").Value = "'2.33"
$ws.Range("E35").Value = "  +10.05%  "
$ws.Range("B36").Value = "InjectiveProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D36").Value = "'38.12"
$ws.Range("E36").Value = "  +0.33%  "
$ws.Range("B37").Value = "OKB"
$ws.Range("C37").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D37").Value = "'51.98"
$ws.Range("E37").Value = "  +0.51%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.0461"
$ws.Range("E38").Value = "  -4.79%  "
$ws.Range("E39").Value = "  +0.25%  "
$ws.Range("D40").Value = "'2.98"
$ws.Range("E40").Value = "  +0.19%  "
$ws.Range("D41").Value = "'3.23"
$ws.Range("E41").Value = "  -8.88%  "
$ws.Range("D42").Value = "'134.89"
$ws.Range("E42").Value = "  +0.38%  "
$ws.Range("D43").Value = "'0.117"
$ws.Range("E43").Value = "  -2.84%  "
$ws.Range("D44").Value = "'1.82"
$ws.Range("E44").Value = "  -3.59%  "
$ws.Range("D45").Value = "'0.269"
$ws.Range("E45").Value = "  -4.39%  "
$ws.Range("B46").Value = "Celestia"
$ws.Range("C46").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D46").Value = "'15.75"
$ws.Range("E46").Value = "  -7.53%  "
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").Value = "'3.65"
$ws.Range("E47").Value = "  -8.33%  "
$ws.Range("D48").Value = "'2.17"
$ws.Range("E48").Value = "  +2.47%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'20.69"
$ws.Range("E49").Value = "  -6.69%  "
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "3.673.93"
$ws.Range("E50").Value = "  +3.49%  "
$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").Value = "2.097.75"
$ws.Range("E51").Value = "  -2.72%  "
